$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Empty out the bold "Youssef Bedoui " run, keeping the (empty) run itself.
# ---------------------------------------------------------------------------
$text = $d.Content.Text
$needle = "Youssef Bedoui "
$idx = $text.IndexOf($needle)
$r = $d.Range($idx, $idx + $needle.Length)
$emptyRunXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:b/></w:rPr><w:t/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($emptyRunXml)

# ---------------------------------------------------------------------------
# 2. The birth-date placeholder "06/10/1998" -> "Test " (stays in the same run).
# ---------------------------------------------------------------------------
$text = $d.Content.Text
$needle = "06/10/1998"
$idx = $text.IndexOf($needle)
$r = $d.Range($idx, $idx + $needle.Length)
$r.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "Test ", 2)

# ---------------------------------------------------------------------------
# 3. The bold CIN-number placeholder "undefined" (right after "CIN N°") -> "Tsemme".
#    Scope the Find to just that run's text so the bold formatting is kept and
#    the other "undefined" placeholders elsewhere in the document are untouched.
# ---------------------------------------------------------------------------
$text = $d.Content.Text
$idx = $text.IndexOf("CIN N°undefined")
$prefix = "CIN N°"
$start = $idx + $prefix.Length
$end = $start + "undefined".Length
$r = $d.Range($start, $end)
$r.Find.Execute("undefined", $true, $false, $false, $false, $false, $true, 1, $false, "Tsemme", 2)

# ---------------------------------------------------------------------------
# 4. "délivrée le undefined" -> "délivrée le Cjtkoy" (only the date placeholder).
# ---------------------------------------------------------------------------
$text = $d.Content.Text
$idx = $text.IndexOf("délivrée le undefined")
$prefix = "délivrée le "
$start = $idx + $prefix.Length
$end = $start + "undefined".Length
$r = $d.Range($start, $end)
$r.Find.Execute("undefined", $true, $false, $false, $false, $false, $true, 1, $false, "Cjtkoy", 2)
